$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Karasjok"
$ws.Range("B14").Value = 69.4448
$ws.Range("C14").Value = 25.3864
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "Approximate location near Karasjok. Taken from https://helikopter.flights/quotes/Price_Taxi_Alta,Norge_Karasjokkommune,Norge.pdf"

$ws.Range("D14").Style = $ws.Range("D9").Style

$ws.Range("E15").Select()
